$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "(~A1+A3)(A1++~A3)(A1+A2)"
$ws.Range("C4").Value = "(A1+A3)(~A1+A2+~A3)"
$ws.Range("C5").Value = "(~A1+A2)"
$ws.Range("C7").Value = "(~A1+A2)(~A1+A3)"

$ws.Range("C9").Select()
